$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-04-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-27 Thursday", 2) | Out-Null
$d.Content.Find.Execute("96×63=6048", $true, $false, $false, $false, $false, $true, 1, $false, "27×55=1485", 2) | Out-Null
$d.Content.Find.Execute("88×16=1408", $true, $false, $false, $false, $false, $true, 1, $false, "86×89=7654", 2) | Out-Null
$d.Content.Find.Execute("46×77=3542", $true, $false, $false, $false, $false, $true, 1, $false, "91×44=4004", 2) | Out-Null
$d.Content.Find.Execute("87×48=4176", $true, $false, $false, $false, $false, $true, 1, $false, "93×42=3906", 2) | Out-Null
$d.Content.Find.Execute("50×40=2000", $true, $false, $false, $false, $false, $true, 1, $false, "91×57=5187", 2) | Out-Null
$d.Content.Find.Execute("84×91=7644", $true, $false, $false, $false, $false, $true, 1, $false, "85×97=8245", 2) | Out-Null
$d.Content.Find.Execute("48×40=1920", $true, $false, $false, $false, $false, $true, 1, $false, "43×100=4300", 2) | Out-Null
$d.Content.Find.Execute("84×45=3780", $true, $false, $false, $false, $false, $true, 1, $false, "12×96=1152", 2) | Out-Null
$d.Content.Find.Execute("27×94=2538", $true, $false, $false, $false, $false, $true, 1, $false, "22×85=1870", 2) | Out-Null
$d.Content.Find.Execute("94×47=4418", $true, $false, $false, $false, $false, $true, 1, $false, "25×67=1675", 2) | Out-Null
$d.Content.Find.Execute("43×77=3311", $true, $false, $false, $false, $false, $true, 1, $false, "99×37=3663", 2) | Out-Null
$d.Content.Find.Execute("11×45=495", $true, $false, $false, $false, $false, $true, 1, $false, "48×62=2976", 2) | Out-Null
$d.Content.Find.Execute("55×51=2805", $true, $false, $false, $false, $false, $true, 1, $false, "87×66=5742", 2) | Out-Null
$d.Content.Find.Execute("54×17=918", $true, $false, $false, $false, $false, $true, 1, $false, "71×59=4189", 2) | Out-Null
$d.Content.Find.Execute("10×85=850", $true, $false, $false, $false, $false, $true, 1, $false, "11×99=1089", 2) | Out-Null
$d.Content.Find.Execute("61×96=5856", $true, $false, $false, $false, $false, $true, 1, $false, "10×23=230", 2) | Out-Null
$d.Content.Find.Execute("50×54=2700", $true, $false, $false, $false, $false, $true, 1, $false, "30×90=2700", 2) | Out-Null
$d.Content.Find.Execute("98×11=1078", $true, $false, $false, $false, $false, $true, 1, $false, "60×79=4740", 2) | Out-Null
$d.Content.Find.Execute("62×46=2852", $true, $false, $false, $false, $false, $true, 1, $false, "55×64=3520", 2) | Out-Null
$d.Content.Find.Execute("10×26=260", $true, $false, $false, $false, $false, $true, 1, $false, "74×84=6216", 2) | Out-Null
$d.Content.Find.Execute("10×29=290", $true, $false, $false, $false, $false, $true, 1, $false, "39×26=1014", 2) | Out-Null
$d.Content.Find.Execute("51×45=2295", $true, $false, $false, $false, $false, $true, 1, $false, "39×33=1287", 2) | Out-Null
$d.Content.Find.Execute("62×33=2046", $true, $false, $false, $false, $false, $true, 1, $false, "33×82=2706", 2) | Out-Null
$d.Content.Find.Execute("29×33=957", $true, $false, $false, $false, $false, $true, 1, $false, "57×50=2850", 2) | Out-Null
$d.Content.Find.Execute("27×22=594", $true, $false, $false, $false, $false, $true, 1, $false, "13×99=1287", 2) | Out-Null
$d.Content.Find.Execute("56×58=3248", $true, $false, $false, $false, $false, $true, 1, $false, "48×30=1440", 2) | Out-Null
$d.Content.Find.Execute("41×60=2460", $true, $false, $false, $false, $false, $true, 1, $false, "11×63=693", 2) | Out-Null
$d.Content.Find.Execute("88×36=3168", $true, $false, $false, $false, $false, $true, 1, $false, "83×63=5229", 2) | Out-Null
$d.Content.Find.Execute("16×53=848", $true, $false, $false, $false, $false, $true, 1, $false, "20×55=1100", 2) | Out-Null
$d.Content.Find.Execute("52×67=3484", $true, $false, $false, $false, $false, $true, 1, $false, "96×11=1056", 2) | Out-Null
$d.Content.Find.Execute("54×64=3456", $true, $false, $false, $false, $false, $true, 1, $false, "94×91=8554", 2) | Out-Null
$d.Content.Find.Execute("35×71=2485", $true, $false, $false, $false, $false, $true, 1, $false, "74×71=5254", 2) | Out-Null
$d.Content.Find.Execute("87×83=7221", $true, $false, $false, $false, $false, $true, 1, $false, "54×73=3942", 2) | Out-Null
$d.Content.Find.Execute("39×20=780", $true, $false, $false, $false, $false, $true, 1, $false, "49×94=4606", 2) | Out-Null
$d.Content.Find.Execute("19×48=912", $true, $false, $false, $false, $false, $true, 1, $false, "36×79=2844", 2) | Out-Null
$d.Content.Find.Execute("12×49=588", $true, $false, $false, $false, $false, $true, 1, $false, "47×68=3196", 2) | Out-Null
$d.Content.Find.Execute("82×58=4756", $true, $false, $false, $false, $false, $true, 1, $false, "97×50=4850", 2) | Out-Null
$d.Content.Find.Execute("11×54=594", $true, $false, $false, $false, $false, $true, 1, $false, "74×75=5550", 2) | Out-Null
$d.Content.Find.Execute("70×11=770", $true, $false, $false, $false, $false, $true, 1, $false, "59×38=2242", 2) | Out-Null
$d.Content.Find.Execute("28×14=392", $true, $false, $false, $false, $false, $true, 1, $false, "71×63=4473", 2) | Out-Null
$d.Content.Find.Execute("73×34=2482", $true, $false, $false, $false, $false, $true, 1, $false, "65×71=4615", 2) | Out-Null
$d.Content.Find.Execute("78×62=4836", $true, $false, $false, $false, $false, $true, 1, $false, "52×15=780", 2) | Out-Null
$d.Content.Find.Execute("32×45=1440", $true, $false, $false, $false, $false, $true, 1, $false, "73×12=876", 2) | Out-Null
$d.Content.Find.Execute("85×27=2295", $true, $false, $false, $false, $false, $true, 1, $false, "80×43=3440", 2) | Out-Null
$d.Content.Find.Execute("81×93=7533", $true, $false, $false, $false, $false, $true, 1, $false, "95×19=1805", 2) | Out-Null
$d.Content.Find.Execute("89×40=3560", $true, $false, $false, $false, $false, $true, 1, $false, "96×48=4608", 2) | Out-Null
$d.Content.Find.Execute("53×44=2332", $true, $false, $false, $false, $false, $true, 1, $false, "12×18=216", 2) | Out-Null
$d.Content.Find.Execute("77×23=1771", $true, $false, $false, $false, $false, $true, 1, $false, "20×32=640", 2) | Out-Null
$d.Content.Find.Execute("77×15=1155", $true, $false, $false, $false, $false, $true, 1, $false, "42×36=1512", 2) | Out-Null
$d.Content.Find.Execute("78×82=6396", $true, $false, $false, $false, $false, $true, 1, $false, "61×91=5551", 2) | Out-Null
$d.Content.Find.Execute("84×69=5796", $true, $false, $false, $false, $false, $true, 1, $false, "38×78=2964", 2) | Out-Null
$d.Content.Find.Execute("38×96=3648", $true, $false, $false, $false, $false, $true, 1, $false, "39×31=1209", 2) | Out-Null
$d.Content.Find.Execute("44×67=2948", $true, $false, $false, $false, $false, $true, 1, $false, "60×58=3480", 2) | Out-Null
$d.Content.Find.Execute("43×88=3784", $true, $false, $false, $false, $false, $true, 1, $false, "67×56=3752", 2) | Out-Null
$d.Content.Find.Execute("83×98=8134", $true, $false, $false, $false, $false, $true, 1, $false, "37×95=3515", 2) | Out-Null
$d.Content.Find.Execute("98×84=8232", $true, $false, $false, $false, $false, $true, 1, $false, "64×77=4928", 2) | Out-Null
$d.Content.Find.Execute("86×74=6364", $true, $false, $false, $false, $false, $true, 1, $false, "71×33=2343", 2) | Out-Null
$d.Content.Find.Execute("92×94=8648", $true, $false, $false, $false, $false, $true, 1, $false, "34×66=2244", 2) | Out-Null
$d.Content.Find.Execute("48×20=960", $true, $false, $false, $false, $false, $true, 1, $false, "30×35=1050", 2) | Out-Null
$d.Content.Find.Execute("17×34=578", $true, $false, $false, $false, $false, $true, 1, $false, "75×65=4875", 2) | Out-Null
$d.Content.Find.Execute("49×59=2891", $true, $false, $false, $false, $false, $true, 1, $false, "13×37=481", 2) | Out-Null
$d.Content.Find.Execute("49×42=2058", $true, $false, $false, $false, $false, $true, 1, $false, "75×37=2775", 2) | Out-Null
$d.Content.Find.Execute("98×18=1764", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=4524", 2) | Out-Null
$d.Content.Find.Execute("13×34=442", $true, $false, $false, $false, $false, $true, 1, $false, "45×100=4500", 2) | Out-Null
$d.Content.Find.Execute("44×10=440", $true, $false, $false, $false, $false, $true, 1, $false, "38×84=3192", 2) | Out-Null
$d.Content.Find.Execute("28×45=1260", $true, $false, $false, $false, $false, $true, 1, $false, "95×93=8835", 2) | Out-Null
$d.Content.Find.Execute("41×11=451", $true, $false, $false, $false, $false, $true, 1, $false, "96×31=2976", 2) | Out-Null
$d.Content.Find.Execute("26×50=1300", $true, $false, $false, $false, $false, $true, 1, $false, "68×26=1768", 2) | Out-Null
$d.Content.Find.Execute("43×40=1720", $true, $false, $false, $false, $false, $true, 1, $false, "77×77=5929", 2) | Out-Null
$d.Content.Find.Execute("96×36=3456", $true, $false, $false, $false, $false, $true, 1, $false, "63×98=6174", 2) | Out-Null
$d.Content.Find.Execute("91×21=1911", $true, $false, $false, $false, $false, $true, 1, $false, "59×85=5015", 2) | Out-Null
$d.Content.Find.Execute("15×11=165", $true, $false, $false, $false, $false, $true, 1, $false, "68×66=4488", 2) | Out-Null
$d.Content.Find.Execute("31×56=1736", $true, $false, $false, $false, $false, $true, 1, $false, "17×80=1360", 2) | Out-Null
$d.Content.Find.Execute("100×42=4200", $true, $false, $false, $false, $false, $true, 1, $false, "11×81=891", 2) | Out-Null
$d.Content.Find.Execute("25×57=1425", $true, $false, $false, $false, $false, $true, 1, $false, "83×70=5810", 2) | Out-Null
$d.Content.Find.Execute("94×15=1410", $true, $false, $false, $false, $false, $true, 1, $false, "87×40=3480", 2) | Out-Null
$d.Content.Find.Execute("61×46=2806", $true, $false, $false, $false, $false, $true, 1, $false, "74×82=6068", 2) | Out-Null
$d.Content.Find.Execute("100×87=8700", $true, $false, $false, $false, $false, $true, 1, $false, "44×58=2552", 2) | Out-Null
$d.Content.Find.Execute("70×82=5740", $true, $false, $false, $false, $false, $true, 1, $false, "26×60=1560", 2) | Out-Null
$d.Content.Find.Execute("51×15=765", $true, $false, $false, $false, $false, $true, 1, $false, "26×91=2366", 2) | Out-Null
$d.Content.Find.Execute("99×65=6435", $true, $false, $false, $false, $false, $true, 1, $false, "23×86=1978", 2) | Out-Null
$d.Content.Find.Execute("21×91=1911", $true, $false, $false, $false, $false, $true, 1, $false, "70×86=6020", 2) | Out-Null
$d.Content.Find.Execute("28×13=364", $true, $false, $false, $false, $false, $true, 1, $false, "41×17=697", 2) | Out-Null
$d.Content.Find.Execute("86×96=8256", $true, $false, $false, $false, $false, $true, 1, $false, "83×36=2988", 2) | Out-Null
$d.Content.Find.Execute("83×81=6723", $true, $false, $false, $false, $false, $true, 1, $false, "19×47=893", 2) | Out-Null
$d.Content.Find.Execute("74×95=7030", $true, $false, $false, $false, $false, $true, 1, $false, "52×57=2964", 2) | Out-Null
$d.Content.Find.Execute("73×85=6205", $true, $false, $false, $false, $false, $true, 1, $false, "76×41=3116", 2) | Out-Null
$d.Content.Find.Execute("52×16=832", $true, $false, $false, $false, $false, $true, 1, $false, "75×48=3600", 2) | Out-Null
$d.Content.Find.Execute("11×73=803", $true, $false, $false, $false, $false, $true, 1, $false, "14×84=1176", 2) | Out-Null
$d.Content.Find.Execute("57×27=1539", $true, $false, $false, $false, $false, $true, 1, $false, "52×64=3328", 2) | Out-Null
$d.Content.Find.Execute("69×60=4140", $true, $false, $false, $false, $false, $true, 1, $false, "71×90=6390", 2) | Out-Null
$d.Content.Find.Execute("91×53=4823", $true, $false, $false, $false, $false, $true, 1, $false, "18×81=1458", 2) | Out-Null
$d.Content.Find.Execute("14×11=154", $true, $false, $false, $false, $false, $true, 1, $false, "50×99=4950", 2) | Out-Null
$d.Content.Find.Execute("62×76=4712", $true, $false, $false, $false, $false, $true, 1, $false, "76×60=4560", 2) | Out-Null
$d.Content.Find.Execute("68×55=3740", $true, $false, $false, $false, $false, $true, 1, $false, "66×44=2904", 2) | Out-Null
$d.Content.Find.Execute("51×37=1887", $true, $false, $false, $false, $false, $true, 1, $false, "75×46=3450", 2) | Out-Null
$d.Content.Find.Execute("90×88=7920", $true, $false, $false, $false, $false, $true, 1, $false, "95×55=5225", 2) | Out-Null
$d.Content.Find.Execute("50×69=3450", $true, $false, $false, $false, $false, $true, 1, $false, "65×80=5200", 2) | Out-Null
$d.Content.Find.Execute("24×21=504", $true, $false, $false, $false, $false, $true, 1, $false, "49×89=4361", 2) | Out-Null
$d.Content.Find.Execute("37×77=2849", $true, $false, $false, $false, $false, $true, 1, $false, "25×75=1875", 2) | Out-Null
